# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1417
$ws1.Range("F3").Value = 2980
$ws1.Range("F5").Value = 276

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1417
$ws4.Range("F3").Value = 2980
$ws4.Range("F6").Value = 276
